$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the date label in A3 (was 2020-07-01_diff -> 2020-06-29_diff)
$ws.Range("A3").Value = "2020-06-29_diff"

# Correct the selection-scope values for that row
$ws.Range("B3").Value = 0.09032787964305289
$ws.Range("C3").Value = 1.422339425567631
$ws.Range("D3").Value = -10.48130760553368
